$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hours on the sprint backlog to be more realistic
$ws.Range("D6").Value = 2
$ws.Range("D9").Value = 2

# Update the active selection to reflect where the author left off editing
$ws.Range("C6").Select() | Out-Null
